$wb = $excel.ActiveWorkbook

# --- ALC row 12 (hunk 0) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 465.4
$ws.Range("I12").Value = 465.4
$ws.Range("K12").Value = 465.4
$ws.Range("M12").Value = -295.4

# --- ALC row 40 (hunk 1) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 45239428
$ws.Range("I40").Value = 35714284
$ws.Range("J40").Value = 50002000
$ws.Range("K40").Value = 35714284
$ws.Range("L40").Value = 50002000
$ws.Range("M40").Value = -35714109
$ws.Range("N40").Value = -50002350

# --- ALC row 80 (hunk 2) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1016.64514
$ws.Range("J80").Value = 1024.1904
$ws.Range("L80").Value = 3072.5712
$ws.Range("N80").Value = -5068.5712

# --- ALC row 83 (hunk 3) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1016.64514
$ws.Range("J83").Value = 1024.1904
$ws.Range("L83").Value = 9217.713599999999
$ws.Range("N83").Value = -19201.7136

# --- ALC row 113 (hunk 4) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 7367.35
$ws.Range("I113").Value = 4144.25
$ws.Range("J113").Value = 9516.083000000001
$ws.Range("K113").Value = 4144.25
$ws.Range("L113").Value = 9516.083000000001
$ws.Range("M113").Value = -890.25
$ws.Range("N113").Value = -16024.083

# --- ARM row 61 (hunk 5) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2471
$ws.Range("I61").Value = 1404.6428
$ws.Range("K61").Value = 1404.6428
$ws.Range("M61").Value = -1192.6428

# --- ARM row 74 (hunk 6) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2792.6667
$ws.Range("J74").Value = 3119.1292
$ws.Range("L74").Value = 3119.1292
$ws.Range("N74").Value = -4867.129199999999

# --- ARM row 76 (hunk 7) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 54444.285
$ws.Range("I76").Value = 38500
$ws.Range("J76").Value = 60822
$ws.Range("K76").Value = 38500
$ws.Range("L76").Value = 60822
$ws.Range("M76").Value = -38162
$ws.Range("N76").Value = -61498

# --- ARM row 77 (hunk 8) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2792.6667
$ws.Range("J77").Value = 3119.1292
$ws.Range("L77").Value = 15595.646
$ws.Range("N77").Value = -24331.646

# --- ARM row 79 (hunk 9) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 54444.285
$ws.Range("I79").Value = 38500
$ws.Range("J79").Value = 60822
$ws.Range("K79").Value = 38500
$ws.Range("L79").Value = 60822
$ws.Range("M79").Value = -37330
$ws.Range("N79").Value = -63162

# --- ARM row 88 (hunk 10) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1544.5151
$ws.Range("I88").Value = 1208.9231
$ws.Range("J88").Value = 1762.65
$ws.Range("K88").Value = 1208.9231
$ws.Range("L88").Value = 1762.65
$ws.Range("M88").Value = -802.9231
$ws.Range("N88").Value = -2574.65

# --- ARM row 91 (hunk 11) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1544.5151
$ws.Range("I91").Value = 1208.9231
$ws.Range("J91").Value = 1762.65
$ws.Range("K91").Value = 1208.9231
$ws.Range("L91").Value = 1762.65
$ws.Range("M91").Value = 195.0769
$ws.Range("N91").Value = -4570.65

# --- ARM row 132 (hunk 12) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5246.2793
$ws.Range("I132").Value = 5407.25
$ws.Range("J132").Value = 4418.4287
$ws.Range("K132").Value = 16221.75
$ws.Range("L132").Value = 13255.2861
$ws.Range("M132").Value = -13691.75
$ws.Range("N132").Value = -18315.2861

# --- ARM row 136 (hunk 13) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2471
$ws.Range("I136").Value = 1404.6428
$ws.Range("K136").Value = 4213.928400000001
$ws.Range("M136").Value = -1663.928400000001

# --- BSM row 82 (hunk 14) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 11084
$ws.Range("I82").Value = 11084
$ws.Range("K82").Value = 11084
$ws.Range("M82").Value = -10701

# --- BSM row 85 (hunk 15) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 11084
$ws.Range("I85").Value = 11084
$ws.Range("K85").Value = 11084
$ws.Range("M85").Value = -9758

# --- BSM row 107 (hunk 16) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1888.8182
$ws.Range("I107").Value = 1819.7368
$ws.Range("K107").Value = 1819.7368
$ws.Range("M107").Value = 100.2632000000001

# --- BSM row 135 (hunk 17) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 77955.266
$ws.Range("J135").Value = 77955.266
$ws.Range("L135").Value = 77955.266
$ws.Range("N135").Value = -88095.266

# --- CRP row 38 (hunk 18) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 13012.667
$ws.Range("I38").Value = 13012.667
$ws.Range("K38").Value = 13012.667
$ws.Range("M38").Value = -12635.667

# --- CRP row 46 (hunk 19) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H46").Value = 13012.667
$ws.Range("I46").Value = 13012.667
$ws.Range("K46").Value = 13012.667
$ws.Range("M46").Value = -12801.667

# --- CRP row 58 (hunk 20) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5392.3335
$ws.Range("I58").Value = 1839
$ws.Range("J58").Value = 12499
$ws.Range("K58").Value = 1839
$ws.Range("L58").Value = 12499
$ws.Range("M58").Value = -1636
$ws.Range("N58").Value = -12905

# --- CRP row 107 (hunk 21) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2620327.5
$ws.Range("I107").Value = 8462177
$ws.Range("J107").Value = 1567.5518
$ws.Range("K107").Value = 8462177
$ws.Range("L107").Value = 1567.5518
$ws.Range("M107").Value = -8460257
$ws.Range("N107").Value = -5407.5518

# --- CRP row 132 (hunk 22) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3315.0908
$ws.Range("I132").Value = 3463.8572
$ws.Range("K132").Value = 10391.5716
$ws.Range("M132").Value = -7861.571599999999

# --- CRP row 134 (hunk 23) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3286.5588
$ws.Range("I134").Value = 3100.7144
$ws.Range("J134").Value = 4153.8335
$ws.Range("K134").Value = 9302.143199999999
$ws.Range("L134").Value = 12461.5005
$ws.Range("M134").Value = -6767.143199999999
$ws.Range("N134").Value = -17531.5005

# --- CRP row 136 (hunk 24) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 5392.3335
$ws.Range("I136").Value = 1839
$ws.Range("J136").Value = 12499
$ws.Range("K136").Value = 5517
$ws.Range("L136").Value = 37497
$ws.Range("M136").Value = -2967
$ws.Range("N136").Value = -42597

# --- CUL row 68 (hunk 25) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 496.66666
$ws.Range("J68").Value = 496.66666
$ws.Range("L68").Value = 1489.99998
$ws.Range("N68").Value = -3111.99998

# --- CUL row 71 (hunk 26) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 496.66666
$ws.Range("J71").Value = 496.66666
$ws.Range("L71").Value = 4469.99994
$ws.Range("N71").Value = -12581.99994

# --- CUL row 92 (hunk 27) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 275.8889
$ws.Range("I92").Value = 252.25
$ws.Range("J92").Value = 294.8
$ws.Range("K92").Value = 756.75
$ws.Range("L92").Value = 884.4000000000001
$ws.Range("M92").Value = 491.25
$ws.Range("N92").Value = -3380.4

# --- CUL row 113 (hunk 28) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1779.4375
$ws.Range("I113").Value = 1880.6666
$ws.Range("K113").Value = 5641.9998
$ws.Range("M113").Value = -3471.9998

# --- CUL row 132 (hunk 29) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1736.5555
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# --- GSM row 80 (hunk 30) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4681.16
$ws.Range("I80").Value = 2188.8462
$ws.Range("J80").Value = 7381.1665
$ws.Range("K80").Value = 2188.8462
$ws.Range("L80").Value = 7381.1665
$ws.Range("M80").Value = -1190.8462
$ws.Range("N80").Value = -9377.166499999999

# --- GSM row 83 (hunk 31) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4681.16
$ws.Range("I83").Value = 2188.8462
$ws.Range("J83").Value = 7381.1665
$ws.Range("K83").Value = 10944.231
$ws.Range("L83").Value = 36905.8325
$ws.Range("M83").Value = -5952.231
$ws.Range("N83").Value = -46889.8325

# --- GSM row 113 (hunk 32) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

# --- GSM row 132 (hunk 33) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3510.3635
$ws.Range("I132").Value = 3482.84
$ws.Range("K132").Value = 10448.52
$ws.Range("M132").Value = -7918.52

# --- LTW row 82 (hunk 34) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2291.1177
$ws.Range("I82").Value = 2249.8
$ws.Range("J82").Value = 2350.1428
$ws.Range("K82").Value = 2249.8
$ws.Range("L82").Value = 2350.1428
$ws.Range("M82").Value = -1888.8
$ws.Range("N82").Value = -3072.1428

# --- LTW row 85 (hunk 35) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2291.1177
$ws.Range("I85").Value = 2249.8
$ws.Range("J85").Value = 2350.1428
$ws.Range("K85").Value = 2249.8
$ws.Range("L85").Value = 2350.1428
$ws.Range("M85").Value = -1001.8
$ws.Range("N85").Value = -4846.1428

# --- LTW row 136 (hunk 36) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3915.1667
$ws.Range("I136").Value = 5330.6665
$ws.Range("K136").Value = 15991.9995
$ws.Range("M136").Value = -13441.9995

# --- WVR row 3 (hunk 37) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 4332.3335
$ws.Range("J3").Value = 3000
$ws.Range("L3").Value = 3000
$ws.Range("N3").Value = -3228

# --- WVR row 61 (hunk 38) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 59700
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

# --- WVR row 86 (hunk 39) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 110162.5
$ws.Range("J86").Value = 110162.5
$ws.Range("L86").Value = 110162.5
$ws.Range("N86").Value = -112408.5

# --- WVR row 89 (hunk 40) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H89").Value = 110162.5
$ws.Range("J89").Value = 110162.5
$ws.Range("L89").Value = 550812.5
$ws.Range("N89").Value = -562044.5

# --- WVR row 107 (hunk 41) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 298.75
$ws.Range("I107").Value = 299.66666
$ws.Range("J107").Value = 292.33334
$ws.Range("K107").Value = 898.9999799999999
$ws.Range("L107").Value = 877.0000200000001
$ws.Range("M107").Value = 1021.00002
$ws.Range("N107").Value = -4717.00002

# --- WVR row 122 (hunk 42) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 17860824
$ws.Range("I122").Value = 26319796
$ws.Range("J122").Value = 2994.3333
$ws.Range("K122").Value = 78959388
$ws.Range("L122").Value = 8982.999899999999
$ws.Range("M122").Value = -78956938
$ws.Range("N122").Value = -13882.9999

# --- WVR row 132 (hunk 43) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1601
$ws.Range("I132").Value = 1345.375
$ws.Range("K132").Value = 4036.125
$ws.Range("M132").Value = -1506.125

# --- WVR row 133 (hunk 44) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 108625
$ws.Range("J133").Value = 108625
$ws.Range("L133").Value = 108625
$ws.Range("N133").Value = -118745
